$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 53 for the November 2016 monthly data, shifting
# everything below (Year to Date / Rolling 12 Months sections) down by one.
$ws.Rows.Item(53).Insert()

# Copy formatting from the October row (now row 52) into the new blank row
# so the new row matches the existing monthly data rows exactly.
$ws.Range("A52:F52").Copy()
$ws.Range("A53:F53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the November 2016 data.
$ws.Range("A53").Value = "November"
$ws.Range("B53").Value = 26875
$ws.Range("C53").Value = 2417
$ws.Range("D53").Value = 21751
$ws.Range("E53").Value = 2181
$ws.Range("F53").Value = 525

# Update the title and "Rolling 12 months" label text for November.
$ws.Range("A2").Value = "by Sector, 2006-November 2016 (Million Cubic Feet)"
$ws.Range("A58").Value = "Rolling 12 Months Ending in November"

# Year to Date figures (now rows 55-57).
$ws.Range("B55").Value = 265108
$ws.Range("C55").Value = 23918
$ws.Range("D55").Value = 211054
$ws.Range("E55").Value = 25547
$ws.Range("F55").Value = 4590

$ws.Range("B56").Value = 258380
$ws.Range("C56").Value = 23094
$ws.Range("D56").Value = 207146
$ws.Range("E56").Value = 23605
$ws.Range("F56").Value = 4534

$ws.Range("B57").Value = 305231
$ws.Range("C57").Value = 27989
$ws.Range("D57").Value = 247622
$ws.Range("E57").Value = 24047
$ws.Range("F57").Value = 5572

# Rolling 12 Months figures (now rows 59-60).
$ws.Range("B59").Value = 280963
$ws.Range("C59").Value = 25171
$ws.Range("D59").Value = 225064
$ws.Range("E59").Value = 25771
$ws.Range("F59").Value = 4957

$ws.Range("B60").Value = 330904
$ws.Range("C60").Value = 30155
$ws.Range("D60").Value = 268501
$ws.Range("E60").Value = 26207
$ws.Range("F60").Value = 6042

Write-Host "Edit complete"
